# "finish login log and attmepts limitation" -- append 4 new daily-entry rows
# (22-25) to the release-history log, mirroring the existing row layout:
#   col A = entry date (date-formatted, style copied from the row above it)
#   col B = entry description (wrap-text style copied from the row above it)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entries: [row, serial-date, description, rowHeight-or-0]
$entries = @(
    @{Row = 22; Date = 44104; Text = "准备linux环境"; Height = 0},
    @{Row = 23; Date = 44113; Text = "在vmware下安装ubuntu 18.04 desktop，user：ritchie/crm2020，administrator： root/crm2020，安装visualstudio，安装python：apt install python3.8，结果安装的是3.6.9；共享文件夹的时候在linux环境里需要使用mount语句"; Height = 36},
    @{Row = 24; Date = 44114; Text = "download运行错误是因为保存文件的文件夹不存在"; Height = 0},
    @{Row = 25; Date = 44116; Text = "login限制 登录大于3次错误限制1小时后登录，大于6次限制1天后登录。两个函数：authentication(user)用来检查是否是新用户，加锁、解锁；userrecrods(user, field)用来记录登录log，生成小时锁和24小时锁"; Height = 36}
)

foreach ($entry in $entries) {
    $r = $entry.Row
    $prev = $r - 1

    # Copy the row-above's formatting (date style on A, wrap-text style on B)
    # down onto the new row before filling in the real values.
    $ws.Range("A" + $prev + ":B" + $prev).Copy()
    $ws.Range("A" + $r + ":B" + $r).PasteSpecial(-4122)

    $ws.Range("A" + $r).Value = $entry.Date
    $ws.Range("B" + $r).Value = $entry.Text

    if ($entry.Height -gt 0) {
        $ws.Rows.Item($r).RowHeight = $entry.Height
    }
}

# Restore the scroll/selection state recorded in the edited workbook.
$ws.Range("H21").Select()
